# Updated cryptos list on Thu Mar  9 17:35:34 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# tracker sheet with the latest scraped figures. Percent cells (E) keep
# their original "  +x.xx%  " / "  -x.xx%  " padded text layout.
#
# Price cells (D) must stay plain text, exactly as scraped (several use a
# '.' thousands separator, e.g. "21.546.37", and trailing zeros like
# "115.50" must be preserved). Whenever the new price string would also
# parse as a plain number (e.g. "288.68"), Excel's automatic type
# detection would otherwise silently convert the cell to a Number and
# round-trip it as a float (losing formatting such as trailing zeros) —
# so the cell is explicitly pre-formatted as Text ("@") before the value
# is assigned, the same way typing into a Text-formatted cell behaves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceCell {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $cell = $ws.Range($CellRef)
    if ($NewValue -match '^[+-]?\d+(\.\d+)?$') {
        # Force Text format first so the numeric-looking string is stored
        # verbatim instead of being auto-converted to a Number.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $NewValue
}

Set-PriceCell "D2" "21.546.37"
$ws.Range("E2").Value = "  -2.51%  "

Set-PriceCell "D3" "1.531.08"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("E5").Value = "  +0.11%  "

Set-PriceCell "D6" "288.68"
$ws.Range("E6").Value = "  -1.20%  "

Set-PriceCell "D7" "0.3875"
$ws.Range("E7").Value = "  -2.50%  "

$ws.Range("E8").Value = "  -2.10%  "

Set-PriceCell "D9" "42.66"
$ws.Range("E9").Value = "  -3.32%  "

$ws.Range("E10").Value = "  -2.40%  "

Set-PriceCell "D11" "1.068"
$ws.Range("E11").Value = "  -1.69%  "

Set-PriceCell "D12" "1.002"
$ws.Range("E12").Value = "  +0.13%  "

Set-PriceCell "D13" "5.723"
$ws.Range("E13").Value = "  +0.51%  "

Set-PriceCell "D14" "18.14"
$ws.Range("E14").Value = "  -4.51%  "

Set-PriceCell "D15" "6.547"
$ws.Range("E15").Value = "  -1.66%  "

Set-PriceCell "D16" "1.537.55"
$ws.Range("E16").Value = "  -1.35%  "

$ws.Range("E17").Value = "  -5.01%  "

Set-PriceCell "D18" "0.06614"
$ws.Range("E18").Value = "  +0.13%  "

Set-PriceCell "D19" "83.52"

$ws.Range("E20").Value = "  +0.18%  "

Set-PriceCell "D21" "6.093"
$ws.Range("E21").Value = "  -3.59%  "

$ws.Range("E22").Value = "  -2.33%  "

Set-PriceCell "D23" "10.79"
$ws.Range("E23").Value = "  -4.45%  "

Set-PriceCell "D24" "2.369"
$ws.Range("E24").Value = "  +0.84%  "

Set-PriceCell "D25" "21.553.80"
$ws.Range("E25").Value = "  -2.53%  "

Set-PriceCell "D26" "2.372"
$ws.Range("E26").Value = "  -3.09%  "

Set-PriceCell "D27" "149.07"
$ws.Range("E27").Value = "  +0.39%  "

Set-PriceCell "D28" "18.31"
$ws.Range("E28").Value = "  -1.76%  "

Set-PriceCell "D29" "4.827"
$ws.Range("E29").Value = "  -0.86%  "

Set-PriceCell "D30" "1.706.11"
$ws.Range("E30").Value = "  -1.54%  "

Set-PriceCell "D31" "116.48"
$ws.Range("E31").Value = "  -2.17%  "

$ws.Range("E32").Value = "  +5.15%  "

Set-PriceCell "D33" "0.9489"
$ws.Range("E33").Value = "  -6.33%  "

Set-PriceCell "D34" "0.07985"
$ws.Range("E34").Value = "  -4.56%  "

Set-PriceCell "D35" "8.517"
$ws.Range("E35").Value = "  -6.14%  "

Set-PriceCell "D36" "5.158"
$ws.Range("E36").Value = "  +0.24%  "

Set-PriceCell "D37" "1.484"
$ws.Range("E37").Value = "  -8.69%  "

Set-PriceCell "D38" "0.02205"
$ws.Range("E38").Value = "  -3.37%  "

Set-PriceCell "D39" "11.28"
$ws.Range("E39").Value = "  +4.50%  "

Set-PriceCell "D40" "0.05883"
$ws.Range("E40").Value = "  -4.27%  "

Set-PriceCell "D41" "0.2019"
$ws.Range("E41").Value = "  -2.31%  "

Set-PriceCell "D42" "1.179"
$ws.Range("E42").Value = "  -3.18%  "

$ws.Range("E43").Value = "  +0.17%  "

Set-PriceCell "D44" "0.5742"
$ws.Range("E44").Value = "  -2.20%  "

Set-PriceCell "D45" "13.22"
$ws.Range("E45").Value = "  +0.51%  "

Set-PriceCell "D46" "3.713"
$ws.Range("E46").Value = "  -1.42%  "

Set-PriceCell "D47" "0.5542"
$ws.Range("E47").Value = "  -1.39%  "

Set-PriceCell "D48" "1.889"
$ws.Range("E48").Value = "  -1.30%  "

Set-PriceCell "D49" "1.155"
$ws.Range("E49").Value = "  +1.24%  "

Set-PriceCell "D50" "115.50"
$ws.Range("E50").Value = "  -3.03%  "

Set-PriceCell "D51" "0.06675"
$ws.Range("E51").Value = "  -2.58%  "
